# Edit script: 
#  1) Add speaker notes ("Michael"/"Alex"/"Nalet") to every slide's Notes page.
#  2) Fix a typo on slide 1: "Michael Utz" -> "Michel Utz".

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Speaker notes per slide (author initials used as placeholder content).
# ---------------------------------------------------------------------------
$notesBySlide = @{
    1 = "Michael"
    2 = "Michael"
    3 = "Alex"
    4 = "Alex"
    5 = "Nalet"
    6 = "Nalet"
    7 = "Nalet"
    8 = "Nalet"
}

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    $notesPage = $slide.NotesPage

    # ppPlaceholderBody = 2. On a notes page only the notes-body placeholder
    # can be created/addressed; AddPlaceholder materializes it (and returns
    # the existing one if it is already there).
    $notesShape = $notesPage.Shapes.AddPlaceholder(2)
    $notesShape.TextFrame.TextRange.Text = $notesBySlide[$i]
}

# ---------------------------------------------------------------------------
# 2) Spelling fix on slide 1: "Nalet Meinen, Alexander Nussbaum, Michael Utz"
#    -> "Nalet Meinen, Alexander Nussbaum, Michel Utz"
# ---------------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
for ($i = 1; $i -le $slide1.Shapes.Count; $i++) {
    $shape = $slide1.Shapes.Item($i)
    if ($shape.HasTextFrame) {
        $tr = $shape.TextFrame.TextRange
        $paraCount = $tr.Paragraphs().Count
        for ($j = 1; $j -le $paraCount; $j++) {
            $para = $tr.Paragraphs($j, 1)
            if ($para.Text -eq "Nalet Meinen, Alexander Nussbaum, Michael Utz") {
                $run = $para.Runs(1, 1)
                $run.Text = "Nalet Meinen, Alexander Nussbaum, Michel Utz"
            }
        }
    }
}
